$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) values - these are stored as text in the sheet,
# so force text entry (avoids Excel auto-converting numeric-looking strings
# like "592.15" into real numbers) and restore the cell's original style
# afterwards so only the value changes.
$priceUpdates = @{
    2  = "66.919.91"
    3  = "2.603.96"
    5  = "592.15"
    6  = "151.92"
    9  = "2.601.81"
    13 = "0.344"
    14 = "27.27"
    15 = "3.078.41"
    17 = "66.789.50"
    18 = "2.602.62"
    19 = "363.05"
    21 = "7.33"
    23 = "2.03"
    24 = "1.00"
    25 = "9.91"
    26 = "66.45"
    27 = "2.737.41"
    29 = "575.58"
    30 = "0.0₂01000"
    32 = "7.68"
    35 = "0.123"
    38 = "156.98"
    39 = "18.93"
    44 = "40.93"
    46 = "16.36"
    47 = "154.36"
    49 = "3.71"
    50 = "21.42"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Updated "Volume(1h)" (column E) values - plain text, never numeric-looking.
$volumeUpdates = @{
    2  = "  -1.93%  "
    3  = "  -0.92%  "
    4  = "  +0.00%  "
    5  = "  -1.31%  "
    6  = "  -3.63%  "
    7  = "  +0.04%  "
    8  = "  +1.48%  "
    9  = "  -0.99%  "
    10 = "  -2.73%  "
    11 = "  +0.32%  "
    12 = "  -2.26%  "
    13 = "  -4.55%  "
    14 = "  -2.37%  "
    15 = "  -0.87%  "
    16 = "  -5.16%  "
    17 = "  -1.19%  "
    18 = "  -0.84%  "
    19 = "  -0.20%  "
    20 = "  -4.49%  "
    21 = "  -6.00%  "
    22 = "  -0.94%  "
    23 = "  -2.94%  "
    24 = "  +0.01%  "
    25 = "  -1.94%  "
    26 = "  -2.62%  "
    27 = "  -0.83%  "
    28 = "  -0.06%  "
    29 = "  -3.80%  "
    30 = "  -5.62%  "
    31 = "  -6.72%  "
    32 = "  -4.76%  "
    33 = "  -3.20%  "
    34 = "  -0.06%  "
    35 = "  -8.73%  "
    36 = "  -4.92%  "
    37 = "  -3.58%  "
    38 = "  +1.56%  "
    39 = "  -3.72%  "
    40 = "  -2.34%  "
    41 = "  -5.15%  "
    42 = "  -6.04%  "
    43 = "  -5.14%  "
    44 = "  -1.50%  "
    46 = "  -0.99%  "
    47 = "  -2.56%  "
    48 = "  -4.41%  "
    49 = "  -1.63%  "
    50 = "  +0.89%  "
    51 = "  -3.35%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
